$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Cells.Item(2, 4).Value = 4989
$ws.Cells.Item(2, 5).Value = 1764
$ws.Cells.Item(2, 6).Value = 1764
$ws.Cells.Item(2, 7).Value = 1688
$ws.Cells.Item(2, 8).Value = 1498
$ws.Cells.Item(2, 9).Value = 1501
$ws.Cells.Item(2, 10).Value = -3
$ws.Cells.Item(2, 11).Value = 27680
$ws.Cells.Item(2, 12).Value = 3048
$ws.Cells.Item(2, 13).Value = 24632
$ws.Cells.Item(2, 14).Value = 24546
$ws.Cells.Item(2, 15).Value = 86
$ws.Cells.Item(2, 16).Value = 291
$ws.Cells.Item(2, 17).Value = 2203
$ws.Cells.Item(2, 18).Value = 1726
$ws.Cells.Item(2, 19).Value = 354
$ws.Cells.Item(2, 20).Value = 128
$ws.Cells.Item(2, 21).Value = 2075
$ws.Cells.Item(2, 22).Value = 3
$ws.Cells.Item(2, 23).Value = 35.37
$ws.Cells.Item(2, 24).Value = 30.03
$ws.Cells.Item(2, 25).Value = 11.41
$ws.Cells.Item(2, 26).Value = 10.04
$ws.Cells.Item(2, 27).Value = 12.37
$ws.Cells.Item(2, 28).Value = 8411.950000000001
$ws.Cells.Item(2, 29).Value = 6116
$ws.Cells.Item(2, 30).Value = 20.21
$ws.Cells.Item(2, 31).Value = 42337
$ws.Cells.Item(2, 32).Value = 2.92
$ws.Cells.Item(2, 33).Value = 173
$ws.Cells.Item(2, 34).Value = 0.14
$ws.Cells.Item(2, 35).Value = 6.68
$ws.Cells.Item(2, 36).Value = 58142204

$ws.Cells.Item(3, 4).Value = 9322
$ws.Cells.Item(3, 5).Value = 886
$ws.Cells.Item(3, 6).Value = 886
$ws.Cells.Item(3, 7).Value = 1095
$ws.Cells.Item(3, 8).Value = 788
$ws.Cells.Item(3, 9).Value = 757
$ws.Cells.Item(3, 10).Value = 31
$ws.Cells.Item(3, 11).Value = 31885
$ws.Cells.Item(3, 12).Value = 6030
$ws.Cells.Item(3, 13).Value = 25855
$ws.Cells.Item(3, 14).Value = 25524
$ws.Cells.Item(3, 15).Value = 331
$ws.Cells.Item(3, 16).Value = 301
$ws.Cells.Item(3, 17).Value = 1622
$ws.Cells.Item(3, 18).Value = -4142
$ws.Cells.Item(3, 19).Value = 1972
$ws.Cells.Item(3, 20).Value = 761
$ws.Cells.Item(3, 21).Value = 861
$ws.Cells.Item(3, 22).Value = 2218
$ws.Cells.Item(3, 23).Value = 9.5
$ws.Cells.Item(3, 24).Value = 8.449999999999999
$ws.Cells.Item(3, 25).Value = 3.02
$ws.Cells.Item(3, 26).Value = 2.65
$ws.Cells.Item(3, 27).Value = 23.32
$ws.Cells.Item(3, 28).Value = 8407.52
$ws.Cells.Item(3, 29).Value = 1269
$ws.Cells.Item(3, 30).Value = 91.23999999999999
$ws.Cells.Item(3, 31).Value = 42476
$ws.Cells.Item(3, 32).Value = 2.73
$ws.Cells.Item(3, 33).Value = 167
$ws.Cells.Item(3, 34).Value = 0.14
$ws.Cells.Item(3, 35).Value = 13.26
$ws.Cells.Item(3, 36).Value = 60096088

$ws.Cells.Item(4, 4).Value = 14642
$ws.Cells.Item(4, 5).Value = 1161
$ws.Cells.Item(4, 6).Value = 1161
$ws.Cells.Item(4, 7).Value = 1003
$ws.Cells.Item(4, 8).Value = 655
$ws.Cells.Item(4, 9).Value = 577
$ws.Cells.Item(4, 10).Value = 78
$ws.Cells.Item(4, 11).Value = 54841
$ws.Cells.Item(4, 12).Value = 17812
$ws.Cells.Item(4, 13).Value = 37029
$ws.Cells.Item(4, 14).Value = 34325
$ws.Cells.Item(4, 15).Value = 2704
$ws.Cells.Item(4, 16).Value = 339
$ws.Cells.Item(4, 17).Value = 3173
$ws.Cells.Item(4, 18).Value = -10000
$ws.Cells.Item(4, 19).Value = 9238
$ws.Cells.Item(4, 20).Value = 810
$ws.Cells.Item(4, 21).Value = 2363
$ws.Cells.Item(4, 22).Value = 9999
$ws.Cells.Item(4, 23).Value = 7.93
$ws.Cells.Item(4, 24).Value = 4.47
$ws.Cells.Item(4, 25).Value = 1.93
$ws.Cells.Item(4, 26).Value = 1.51
$ws.Cells.Item(4, 27).Value = 48.1
$ws.Cells.Item(4, 28).Value = 10068.46
$ws.Cells.Item(4, 29).Value = 874
$ws.Cells.Item(4, 30).Value = 88.15000000000001
$ws.Cells.Item(4, 31).Value = 50769
$ws.Cells.Item(4, 32).Value = 1.52
$ws.Cells.Item(4, 33).Value = 148
$ws.Cells.Item(4, 34).Value = 0.19
$ws.Cells.Item(4, 35).Value = 17.35
$ws.Cells.Item(4, 36).Value = 67615715

$ws.Cells.Item(5, 4).Value = 19723
$ws.Cells.Item(5, 5).Value = 1654
$ws.Cells.Item(5, 6).Value = 1654
$ws.Cells.Item(5, 7).Value = 1533
$ws.Cells.Item(5, 8).Value = 1251
$ws.Cells.Item(5, 9).Value = 1086
$ws.Cells.Item(5, 10).Value = 165
$ws.Cells.Item(5, 11).Value = 63494
$ws.Cells.Item(5, 12).Value = 18865
$ws.Cells.Item(5, 13).Value = 44629
$ws.Cells.Item(5, 14).Value = 40291
$ws.Cells.Item(5, 15).Value = 4787
$ws.Cells.Item(5, 16).Value = 340
$ws.Cells.Item(5, 17).Value = 3719
$ws.Cells.Item(5, 18).Value = -3547
$ws.Cells.Item(5, 19).Value = 4676
$ws.Cells.Item(5, 20).Value = 752
$ws.Cells.Item(5, 21).Value = 2968
$ws.Cells.Item(5, 22).Value = 7952
$ws.Cells.Item(5, 23).Value = 8.380000000000001
$ws.Cells.Item(5, 24).Value = 6.34
$ws.Cells.Item(5, 25).Value = 2.91
$ws.Cells.Item(5, 26).Value = 2.11
$ws.Cells.Item(5, 27).Value = 42.27
$ws.Cells.Item(5, 28).Value = 11781.26
$ws.Cells.Item(5, 29).Value = 1602
$ws.Cells.Item(5, 30).Value = 85.5
$ws.Cells.Item(5, 31).Value = 59336
$ws.Cells.Item(5, 32).Value = 2.31
$ws.Cells.Item(5, 33).Value = 148
$ws.Cells.Item(5, 34).Value = 0.11
$ws.Cells.Item(5, 35).Value = 9.25
$ws.Cells.Item(5, 36).Value = 67908527

$ws.Cells.Item(6, 4).Value = 24170
$ws.Cells.Item(6, 5).Value = 729
$ws.Cells.Item(6, 6).Value = 729
$ws.Cells.Item(6, 7).Value = 1307
$ws.Cells.Item(6, 8).Value = 159
$ws.Cells.Item(6, 9).Value = 479
$ws.Cells.Item(6, 11).Value = 79595
$ws.Cells.Item(6, 12).Value = 23324
$ws.Cells.Item(6, 13).Value = 56272
$ws.Cells.Item(6, 14).Value = 51369
$ws.Cells.Item(6, 16).Value = 417
$ws.Cells.Item(6, 17).Value = 4915
$ws.Cells.Item(6, 18).Value = -12607
$ws.Cells.Item(6, 19).Value = 8905
$ws.Cells.Item(6, 20).Value = 972
$ws.Cells.Item(6, 21).Value = 3943
$ws.Cells.Item(6, 22).Value = 6545
$ws.Cells.Item(6, 23).Value = 3.02
$ws.Cells.Item(6, 24).Value = 0.66
$ws.Cells.Item(6, 25).Value = 1.04
$ws.Cells.Item(6, 26).Value = 0.22
$ws.Cells.Item(6, 27).Value = 41.45
$ws.Cells.Item(6, 28).Value = 12219.62
$ws.Cells.Item(6, 29).Value = 613
$ws.Cells.Item(6, 30).Value = 168.01
$ws.Cells.Item(6, 31).Value = 64897
$ws.Cells.Item(6, 32).Value = 1.59
$ws.Cells.Item(6, 33).Value = 127
$ws.Cells.Item(6, 34).Value = 0.12
$ws.Cells.Item(6, 35).Value = 20.99
$ws.Cells.Item(6, 36).Value = 83387773

$ws.Cells.Item(7, 4).Value = 30805
$ws.Cells.Item(7, 5).Value = 1962
$ws.Cells.Item(7, 7).Value = 2672
$ws.Cells.Item(7, 8).Value = 1543
$ws.Cells.Item(7, 9).Value = 1755
$ws.Cells.Item(7, 11).Value = 85983
$ws.Cells.Item(7, 12).Value = 27544
$ws.Cells.Item(7, 13).Value = 58439
$ws.Cells.Item(7, 14).Value = 53758
$ws.Cells.Item(7, 16).Value = 423
$ws.Cells.Item(7, 17).Value = 6911
$ws.Cells.Item(7, 18).Value = -2647
$ws.Cells.Item(7, 19).Value = 20
$ws.Cells.Item(7, 20).Value = 988
$ws.Cells.Item(7, 21).Value = 4463
$ws.Cells.Item(7, 23).Value = 6.37
$ws.Cells.Item(7, 24).Value = 5.01
$ws.Cells.Item(7, 25).Value = 3.34
$ws.Cells.Item(7, 26).Value = 1.86
$ws.Cells.Item(7, 27).Value = 47.13
$ws.Cells.Item(7, 29).Value = 2090
$ws.Cells.Item(7, 30).Value = 76.06999999999999
$ws.Cells.Item(7, 31).Value = 65567
$ws.Cells.Item(7, 32).Value = 2.42
$ws.Cells.Item(7, 33).Value = 144
$ws.Cells.Item(7, 34).Value = 0.09
$ws.Cells.Item(7, 35).Value = 7.09

$ws.Cells.Item(8, 4).Value = 37690
$ws.Cells.Item(8, 5).Value = 4040
$ws.Cells.Item(8, 7).Value = 4839
$ws.Cells.Item(8, 8).Value = 3176
$ws.Cells.Item(8, 9).Value = 3333
$ws.Cells.Item(8, 11).Value = 92500
$ws.Cells.Item(8, 12).Value = 30843
$ws.Cells.Item(8, 13).Value = 61657
$ws.Cells.Item(8, 14).Value = 56871
$ws.Cells.Item(8, 16).Value = 423
$ws.Cells.Item(8, 17).Value = 7175
$ws.Cells.Item(8, 18).Value = -3713
$ws.Cells.Item(8, 19).Value = 638
$ws.Cells.Item(8, 20).Value = 1089
$ws.Cells.Item(8, 21).Value = 5676
$ws.Cells.Item(8, 23).Value = 10.72
$ws.Cells.Item(8, 24).Value = 8.43
$ws.Cells.Item(8, 25).Value = 6.02
$ws.Cells.Item(8, 26).Value = 3.56
$ws.Cells.Item(8, 27).Value = 50.02
$ws.Cells.Item(8, 29).Value = 3861
$ws.Cells.Item(8, 30).Value = 41.18
$ws.Cells.Item(8, 31).Value = 69285
$ws.Cells.Item(8, 32).Value = 2.29
$ws.Cells.Item(8, 33).Value = 154
$ws.Cells.Item(8, 34).Value = 0.1
$ws.Cells.Item(8, 35).Value = 3.98

$ws.Cells.Item(9, 4).Value = 43696
$ws.Cells.Item(9, 5).Value = 5428
$ws.Cells.Item(9, 7).Value = 6314
$ws.Cells.Item(9, 8).Value = 4237
$ws.Cells.Item(9, 9).Value = 4415
$ws.Cells.Item(9, 11).Value = 99710
$ws.Cells.Item(9, 12).Value = 33834
$ws.Cells.Item(9, 13).Value = 65876
$ws.Cells.Item(9, 14).Value = 61064
$ws.Cells.Item(9, 16).Value = 423
$ws.Cells.Item(9, 17).Value = 8254
$ws.Cells.Item(9, 18).Value = -3857
$ws.Cells.Item(9, 19).Value = 690
$ws.Cells.Item(9, 20).Value = 1189
$ws.Cells.Item(9, 21).Value = 6348
$ws.Cells.Item(9, 23).Value = 12.42
$ws.Cells.Item(9, 24).Value = 9.699999999999999
$ws.Cells.Item(9, 25).Value = 7.49
$ws.Cells.Item(9, 26).Value = 4.41
$ws.Cells.Item(9, 27).Value = 51.36
$ws.Cells.Item(9, 29).Value = 5115
$ws.Cells.Item(9, 30).Value = 31.09
$ws.Cells.Item(9, 31).Value = 74394
$ws.Cells.Item(9, 32).Value = 2.14
$ws.Cells.Item(9, 33).Value = 167
$ws.Cells.Item(9, 34).Value = 3.27
